$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1700
$ws.Range("J32").Value = 1700
$ws.Range("L32").Value = 1700
$ws.Range("N32").Value = -2352
$ws.Range("H94").Value = 4750
$ws.Range("I94").Value = 4750
$ws.Range("K94").Value = 4750
$ws.Range("M94").Value = -4299
$ws.Range("H112").Value = 1305.5555
$ws.Range("J112").Value = 1450
$ws.Range("L112").Value = 4350
$ws.Range("N112").Value = -6566
$ws.Range("H113").Value = 10742.857
$ws.Range("I113").Value = 2266.6667
$ws.Range("J113").Value = 26000
$ws.Range("K113").Value = 2266.6667
$ws.Range("L113").Value = 26000
$ws.Range("M113").Value = 987.3332999999998
$ws.Range("N113").Value = -32508
$ws.Range("H114").Value = 37738.5
$ws.Range("J114").Value = 37738.5
$ws.Range("L114").Value = 37738.5
$ws.Range("N114").Value = -46416.5
$ws.Range("H115").Value = 839.06665
$ws.Range("I115").Value = 458.6
$ws.Range("K115").Value = 1375.8
$ws.Range("M115").Value = 191.1999999999998
$ws.Range("H132").Value = 752.4483
$ws.Range("I132").Value = 533.9583
$ws.Range("J132").Value = 1801.2
$ws.Range("K132").Value = 1601.8749
$ws.Range("L132").Value = 5403.6
$ws.Range("M132").Value = 928.1251
$ws.Range("N132").Value = -10463.6
$ws.Range("H137").Value = 1502.4828
$ws.Range("I137").Value = 1732.5625
$ws.Range("J137").Value = 1219.3077
$ws.Range("K137").Value = 5197.6875
$ws.Range("L137").Value = 3657.9231
$ws.Range("M137").Value = -2647.6875
$ws.Range("N137").Value = -8757.9231
$ws.Range("H138").Value = 20836348
$ws.Range("I138").Value = 1841.5
$ws.Range("J138").Value = 38465544
$ws.Range("K138").Value = 5524.5
$ws.Range("L138").Value = 115396632
$ws.Range("M138").Value = -384.5
$ws.Range("N138").Value = -115406912
$ws.Range("H141").Value = 3261.6428
$ws.Range("I141").Value = 2755.8
$ws.Range("K141").Value = 8267.400000000001
$ws.Range("M141").Value = -3087.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 11969.25
$ws.Range("I43").Value = 8700
$ws.Range("K43").Value = 8700
$ws.Range("M43").Value = -8387
$ws.Range("H132").Value = 2114.279
$ws.Range("I132").Value = 831.10345
$ws.Range("J132").Value = 4772.2856
$ws.Range("K132").Value = 2493.31035
$ws.Range("L132").Value = 14316.8568
$ws.Range("M132").Value = 36.68965000000026
$ws.Range("N132").Value = -19376.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1257.2069
$ws.Range("I99").Value = 1136.1428
$ws.Range("J99").Value = 1575
$ws.Range("K99").Value = 1136.1428
$ws.Range("L99").Value = 1575
$ws.Range("M99").Value = 361.8571999999999
$ws.Range("N99").Value = -4571
$ws.Range("H107").Value = 573.4167
$ws.Range("I107").Value = 534.63635
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 534.63635
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1385.36365
$ws.Range("N107").Value = -4840
$ws.Range("H134").Value = 2116.1428
$ws.Range("I134").Value = 961.6667
$ws.Range("J134").Value = 5002.3335
$ws.Range("K134").Value = 2885.0001
$ws.Range("L134").Value = 15007.0005
$ws.Range("M134").Value = -350.0001000000002
$ws.Range("N134").Value = -20077.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1781.1031
$ws.Range("I31").Value = 663.6531
$ws.Range("J31").Value = 2921.8333
$ws.Range("K31").Value = 663.6531
$ws.Range("L31").Value = 2921.8333
$ws.Range("M31").Value = -368.6531
$ws.Range("N31").Value = -3511.8333
$ws.Range("H34").Value = 1781.1031
$ws.Range("I34").Value = 663.6531
$ws.Range("J34").Value = 2921.8333
$ws.Range("K34").Value = 663.6531
$ws.Range("L34").Value = 2921.8333
$ws.Range("M34").Value = -461.6531
$ws.Range("N34").Value = -3325.8333
$ws.Range("H58").Value = 2198.2222
$ws.Range("I58").Value = 2560.8
$ws.Range("J58").Value = 1745
$ws.Range("K58").Value = 2560.8
$ws.Range("L58").Value = 1745
$ws.Range("M58").Value = -2357.8
$ws.Range("N58").Value = -2151
$ws.Range("H62").Value = 4368.9473
$ws.Range("I62").Value = 4121
$ws.Range("J62").Value = 4457.5
$ws.Range("K62").Value = 4121
$ws.Range("L62").Value = 4457.5
$ws.Range("M62").Value = -3497
$ws.Range("N62").Value = -5705.5
$ws.Range("H65").Value = 4368.9473
$ws.Range("I65").Value = 4121
$ws.Range("J65").Value = 4457.5
$ws.Range("K65").Value = 20605
$ws.Range("L65").Value = 22287.5
$ws.Range("M65").Value = -17485
$ws.Range("N65").Value = -28527.5
$ws.Range("H105").Value = 79240
$ws.Range("I105").Value = 102563
$ws.Range("K105").Value = 102563
$ws.Range("M105").Value = -100816
$ws.Range("H122").Value = 1107.2142
$ws.Range("I122").Value = 975.0833
$ws.Range("K122").Value = 2925.2499
$ws.Range("M122").Value = -475.2498999999998
$ws.Range("H132").Value = 2631.5417
$ws.Range("I132").Value = 2051.2144
$ws.Range("J132").Value = 3444
$ws.Range("K132").Value = 6153.6432
$ws.Range("L132").Value = 10332
$ws.Range("M132").Value = -3623.6432
$ws.Range("N132").Value = -15392
$ws.Range("H134").Value = 3113.9348
$ws.Range("I134").Value = 4185.963
$ws.Range("J134").Value = 1590.5264
$ws.Range("K134").Value = 12557.889
$ws.Range("L134").Value = 4771.5792
$ws.Range("M134").Value = -10022.889
$ws.Range("N134").Value = -9841.5792
$ws.Range("H136").Value = 2198.2222
$ws.Range("I136").Value = 2560.8
$ws.Range("J136").Value = 1745
$ws.Range("K136").Value = 7682.400000000001
$ws.Range("L136").Value = 5235
$ws.Range("M136").Value = -5132.400000000001
$ws.Range("N136").Value = -10335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 255000
$ws.Range("J9").Value = 255000
$ws.Range("L9").Value = 765000
$ws.Range("N9").Value = -765448
$ws.Range("H12").Value = 61.85
$ws.Range("J12").Value = 86.35714
$ws.Range("L12").Value = 259.07142
$ws.Range("N12").Value = -605.07142
$ws.Range("H68").Value = 1920542.1
$ws.Range("I68").Value = 3906146.5
$ws.Range("J68").Value = 1124.4667
$ws.Range("K68").Value = 11718439.5
$ws.Range("L68").Value = 3373.4001
$ws.Range("M68").Value = -11717628.5
$ws.Range("N68").Value = -4995.4001
$ws.Range("H71").Value = 1920542.1
$ws.Range("I71").Value = 3906146.5
$ws.Range("J71").Value = 1124.4667
$ws.Range("K71").Value = 35155318.5
$ws.Range("L71").Value = 10120.2003
$ws.Range("M71").Value = -35151262.5
$ws.Range("N71").Value = -18232.2003
$ws.Range("H75").Value = 2784.3635
$ws.Range("I75").Value = 671
$ws.Range("J75").Value = 3576.875
$ws.Range("K75").Value = 2013
$ws.Range("L75").Value = 10730.625
$ws.Range("M75").Value = -1015
$ws.Range("N75").Value = -12726.625
$ws.Range("H78").Value = 2784.3635
$ws.Range("I78").Value = 671
$ws.Range("J78").Value = 3576.875
$ws.Range("K78").Value = 6039
$ws.Range("L78").Value = 32191.875
$ws.Range("M78").Value = -1047
$ws.Range("N78").Value = -42175.875
$ws.Range("H82").Value = 149064.86
$ws.Range("J82").Value = 149064.86
$ws.Range("L82").Value = 447194.58
$ws.Range("N82").Value = -448006.58
$ws.Range("H85").Value = 149064.86
$ws.Range("J85").Value = 149064.86
$ws.Range("L85").Value = 447194.58
$ws.Range("N85").Value = -450002.58
$ws.Range("H107").Value = 817.5238000000001
$ws.Range("J107").Value = 712.6
$ws.Range("L107").Value = 2137.8
$ws.Range("N107").Value = -5977.8
$ws.Range("H131").Value = 1295.1515
$ws.Range("J131").Value = 1326.1702
$ws.Range("L131").Value = 3978.5106
$ws.Range("N131").Value = -14058.5106

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H102").Value = 1145.8182
$ws.Range("I102").Value = 1132.2222
$ws.Range("K102").Value = 1132.2222
$ws.Range("M102").Value = 489.7778000000001
$ws.Range("H122").Value = 2174.8108
$ws.Range("I122").Value = 1930.6428
$ws.Range("J122").Value = 2934.4443
$ws.Range("K122").Value = 5791.928400000001
$ws.Range("L122").Value = 8803.332900000001
$ws.Range("M122").Value = -3341.928400000001
$ws.Range("N122").Value = -13703.3329
$ws.Range("H132").Value = 2716.111
$ws.Range("I132").Value = 1570.72
$ws.Range("J132").Value = 5319.273
$ws.Range("K132").Value = 4712.16
$ws.Range("L132").Value = 15957.819
$ws.Range("M132").Value = -2182.16
$ws.Range("N132").Value = -21017.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 18520398
$ws.Range("I122").Value = 37038404
$ws.Range("J122").Value = 2393.3333
$ws.Range("K122").Value = 111115212
$ws.Range("L122").Value = 7179.999899999999
$ws.Range("M122").Value = -111112762
$ws.Range("N122").Value = -12079.9999
$ws.Range("H136").Value = 7578230.5
$ws.Range("I136").Value = 1949.6923
$ws.Range("J136").Value = 10755380
$ws.Range("K136").Value = 5849.0769
$ws.Range("L136").Value = 32266140
$ws.Range("M136").Value = -3299.0769
$ws.Range("N136").Value = -32271240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 1676836.5
$ws.Range("J45").Value = 12203.8
$ws.Range("L45").Value = 12203.8
$ws.Range("N45").Value = -13185.8
